$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n MATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight`nWHERE diag.primary_disease_site IN ['Mouth']`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,`n       coalesce(demo.sex, '') AS Sex,`n       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`ncoalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS ``Weight (kg)``,`n       coalesce(diag.best_response, '') AS ``Response to Treatment``,`n       coalesce(co.cohort_description, '') AS ``Cohort```norder by c.case_id asc`nlimit 100"

$ws.Range("B2").Value = $newCasesQuery

$ws.Range("C3").Select()
